$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Ntn1"
$ws.Range("C2").Value = "Unc5c"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 1.644726333333334
$ws.Range("H2").Value = 4.934179
$ws.Range("I2").Value = 0.03084360558270512
$ws.Range("J2").Value = 0.03084360558270512
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.009159
$ws.Range("N2").Value = 0.027477
$ws.Range("O2").Value = 0.01850291816413234
$ws.Range("P2").Value = 0.01850291816413234
$ws.Range("Q2").Value = 0.015064048487
$ws.Range("R2").Value = 0.135576436383
$ws.Range("S2").Value = 0.0005706967099835682
$ws.Range("T2").Value = 0.0005706967099835681

$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Ntn1"
$ws.Range("C3").Value = "Unc5c"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 1.644726333333334
$ws.Range("H3").Value = 4.934179
$ws.Range("I3").Value = 0.03084360558270512
$ws.Range("J3").Value = 0.03084360558270512
$ws.Range("K3").Value = 2
$ws.Range("L3").Value = 0.6666666666666666
$ws.Range("M3").Value = 0.1504516666666667
$ws.Range("N3").Value = 0.451355
$ws.Range("O3").Value = 0.3039409188765859
$ws.Range("P3").Value = 0.3039409188765859
$ws.Range("Q3").Value = 0.2474518180605556
$ws.Range("R3").Value = 2.227066362545
$ws.Range("S3").Value = 0.009374633822274389
$ws.Range("T3").Value = 0.009374633822274388

$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Ntn1"
$ws.Range("C4").Value = "Unc5c"
$ws.Range("D4").Value = "MuSCs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 1.644726333333334
$ws.Range("H4").Value = 4.934179
$ws.Range("I4").Value = 0.03084360558270512
$ws.Range("J4").Value = 0.03084360558270512
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.3353923333333333
$ws.Range("N4").Value = 1.006177
$ws.Range("O4").Value = 0.6775561629592817
$ws.Range("P4").Value = 0.6775561629592817
$ws.Range("Q4").Value = 0.5516286026314445
$ws.Range("R4").Value = 4.964657423683001
$ws.Range("S4").Value = 0.02089827505044716
$ws.Range("T4").Value = 0.02089827505044716

$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Ntn1"
$ws.Range("C5").Value = "Unc5c"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 44.154177
$ws.Range("H5").Value = 132.462531
$ws.Range("I5").Value = 0.828024694817689
$ws.Range("J5").Value = 0.828024694817689
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.009159
$ws.Range("N5").Value = 0.027477
$ws.Range("O5").Value = 0.01850291816413234
$ws.Range("P5").Value = 0.01850291816413234
$ws.Range("Q5").Value = 0.4044081071430001
$ws.Range("R5").Value = 3.639672964287
$ws.Range("S5").Value = 0.01532087316609236
$ws.Range("T5").Value = 0.01532087316609235

$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Ntn1"
$ws.Range("C6").Value = "Unc5c"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 44.154177
$ws.Range("H6").Value = 132.462531
$ws.Range("I6").Value = 0.828024694817689
$ws.Range("J6").Value = 0.828024694817689
$ws.Range("K6").Value = 2
$ws.Range("L6").Value = 0.6666666666666666
$ws.Range("M6").Value = 0.1504516666666667
$ws.Range("N6").Value = 0.451355
$ws.Range("O6").Value = 0.3039409188765859
$ws.Range("P6").Value = 0.3039409188765859
$ws.Range("Q6").Value = 6.643069519945001
$ws.Range("R6").Value = 59.78762567950501
$ws.Range("S6").Value = 0.251670586595393
$ws.Range("T6").Value = 0.251670586595393

$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Ntn1"
$ws.Range("C7").Value = "Unc5c"
$ws.Range("D7").Value = "MuSCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 44.154177
$ws.Range("H7").Value = 132.462531
$ws.Range("I7").Value = 0.828024694817689
$ws.Range("J7").Value = 0.828024694817689
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 0.3353923333333333
$ws.Range("N7").Value = 1.006177
$ws.Range("O7").Value = 0.6775561629592817
$ws.Range("P7").Value = 0.6775561629592817
$ws.Range("Q7").Value = 14.808972450443
$ws.Range("R7").Value = 133.280752053987
$ws.Range("S7").Value = 0.5610332350562036
$ws.Range("T7").Value = 0.5610332350562036

$ws.Range("A8").Value = "MuSCs"
$ws.Range("B8").Value = "Ntn1"
$ws.Range("C8").Value = "Unc5c"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 7.525807
$ws.Range("H8").Value = 22.577421
$ws.Range("I8").Value = 0.1411316995996059
$ws.Range("J8").Value = 0.1411316995996059
$ws.Range("K8").Value = 2
$ws.Range("L8").Value = 0.6666666666666666
$ws.Range("M8").Value = 0.009159
$ws.Range("N8").Value = 0.027477
$ws.Range("O8").Value = 0.01850291816413234
$ws.Range("P8").Value = 0.01850291816413234
$ws.Range("Q8").Value = 0.06892886631300001
$ws.Range("R8").Value = 0.620359796817
$ws.Range("S8").Value = 0.002611348288056417
$ws.Range("T8").Value = 0.002611348288056416

$ws.Range("A9").Value = "MuSCs"
$ws.Range("B9").Value = "Ntn1"
$ws.Range("C9").Value = "Unc5c"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 7.525807
$ws.Range("H9").Value = 22.577421
$ws.Range("I9").Value = 0.1411316995996059
$ws.Range("J9").Value = 0.1411316995996059
$ws.Range("K9").Value = 2
$ws.Range("L9").Value = 0.6666666666666666
$ws.Range("M9").Value = 0.1504516666666667
$ws.Range("N9").Value = 0.451355
$ws.Range("O9").Value = 0.3039409188765859
$ws.Range("P9").Value = 0.3039409188765859
$ws.Range("Q9").Value = 1.132270206161667
$ws.Range("R9").Value = 10.190431855455
$ws.Range("S9").Value = 0.04289569845891851
$ws.Range("T9").Value = 0.0428956984589185

$ws.Range("A10").Value = "MuSCs"
$ws.Range("B10").Value = "Ntn1"
$ws.Range("C10").Value = "Unc5c"
$ws.Range("D10").Value = "MuSCs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 7.525807
$ws.Range("H10").Value = 22.577421
$ws.Range("I10").Value = 0.1411316995996059
$ws.Range("J10").Value = 0.1411316995996059
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 0.3353923333333333
$ws.Range("N10").Value = 1.006177
$ws.Range("O10").Value = 0.6775561629592817
$ws.Range("P10").Value = 0.6775561629592817
$ws.Range("Q10").Value = 2.524097969946334
$ws.Range("R10").Value = 22.716881729517
$ws.Range("S10").Value = 0.09562465285263096
$ws.Range("T10").Value = 0.09562465285263096
